$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Handoff transform failed" -> "Ready for handoff"
#    (shared across Overview!B2/C2 and zh-cn!B2 / de-de!B2)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = "Ready for handoff"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: a handoff has actually gone out for the .md file (row 2) -
#    record the generated xlf handoff file (with hyperlink), the handoff
#    timestamp, and flip the handoff reason from "Ignored" to "Include".
# ---------------------------------------------------------------------------
$zhFileName = "c6618b96-2854-42be-a23b-85bdba8859bd.e5df2e3a7d9ef8081c021184299b731427242f3f.zh-cn.xlf"
$zhUrl = "https://github.com/OpenLocalizationTest/oltest/blob/dd6492fda01ca6011008b594d85bc97bf16184d0/e2e/$zhFileName"

$wsZhCn.Range("C2").Value = $zhFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhUrl, "", "", $zhFileName)
$wsZhCn.Range("C2").Font.Underline = $true
$wsZhCn.Range("C2").Font.Color = 15570276

$wsZhCn.Range("D2").Value = "2016-02-15 04:10:02"
$wsZhCn.Range("H2").Value = "Include"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment, with the de-de handoff file/time.
# ---------------------------------------------------------------------------
$deFileName = "c6618b96-2854-42be-a23b-85bdba8859bd.e5df2e3a7d9ef8081c021184299b731427242f3f.de-de.xlf"
$deUrl = "https://github.com/OpenLocalizationTest/oltest/blob/dd6492fda01ca6011008b594d85bc97bf16184d0/e2e/$deFileName"

$wsDeDe.Range("C2").Value = $deFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $deUrl, "", "", $deFileName)
$wsDeDe.Range("C2").Font.Underline = $true
$wsDeDe.Range("C2").Font.Color = 15570276

$wsDeDe.Range("D2").Value = "2016-02-15 04:10:16"
$wsDeDe.Range("H2").Value = "Include"
